$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.601.52'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.363.71'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.62'
$ws.Range("E5").Value = '  +5.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.87'
$ws.Range("E6").Value = '  -7.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.09'
$ws.Range("E10").Value = '  -7.93%  '
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.42'
$ws.Range("E12").Value = '  -6.80%  '
$ws.Range("E13").Value = '  -5.28%  '
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.30'
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.717.41'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.364.60'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.558.15'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("E19").Value = '  +5.62%  '
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.24'
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("E22").Value = '  +6.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '275.99'
$ws.Range("E23").Value = '  +10.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.30'
$ws.Range("E24").Value = '  -8.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.79'
$ws.Range("E25").Value = '  +8.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -4.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.77'
$ws.Range("E28").Value = '  +5.72%  '
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.87'
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.10'
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '35.54'
$ws.Range("E32").Value = '  -8.20%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0900'
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.97'
$ws.Range("E34").Value = '  +2.84%  '
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.62'
$ws.Range("E36").Value = '  -7.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0358'
$ws.Range("E37").Value = '  -5.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.89'
$ws.Range("E38").Value = '  -6.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.87'
$ws.Range("E39").Value = '  +4.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("E41").Value = '  +1.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.228'
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '68.79'
$ws.Range("E43").Value = '  -4.80%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.98'
$ws.Range("E45").Value = '  +40.08%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '115.40'
$ws.Range("E46").Value = '  +4.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.92'
$ws.Range("E47").Value = '  -5.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.45'
$ws.Range("E48").Value = '  -4.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.02'
$ws.Range("E49").Value = '  -2.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.588.71'
$ws.Range("E50").Value = '  +6.47%  '
$ws.Range("E51").Value = '  -2.88%  '
